$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 844; this shifts the existing rows 844:956
# down to 845:957 (carrying their data and formatting with them), matching
# the diff's observed cascading shift of every record one row down.
$ws.Rows.Item(844).Insert()

# Populate the newly-inserted row 844 with the new weekly record.
$ws.Cells.Item(844, 1).Value = 8
$ws.Cells.Item(844, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(844, 3).Value = "Coquimbo"
$ws.Cells.Item(844, 4).Value = Get-Date -Year 2023 -Month 7 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(844, 5).Value = 4
$ws.Cells.Item(844, 6).Value = 100112045
$ws.Cells.Item(844, 7).Value = "Zapallo"
$ws.Cells.Item(844, 8).Value = "Camote"
$ws.Cells.Item(844, 9).Value = "1a (guarda)"
$ws.Cells.Item(844, 10).Value = 1800
$ws.Cells.Item(844, 11).Value = 650
$ws.Cells.Item(844, 12).Value = 700
$ws.Cells.Item(844, 13).Value = 675
$ws.Cells.Item(844, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(844, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(844, 16).Value = 675
$ws.Cells.Item(844, 17).Value = 1
$ws.Cells.Item(844, 18).Value = "Hortaliza"
